$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived LR-pair data: rows 2-10, columns A-T.
# Columns A-D are cluster/gene labels (text); E-T are numeric metrics.
$data = @(
        @("ECs","Ccl11","Ackr4","ECs","2","0.6666666666666666","0.259924","0.779772","0.0006491957374851489","0.0006491957374851488","1","0.3333333333333333","0.07058066666666667","0.211742","0.03405243646146196","0.03405243646146196","0.01834560920266667","0.165110482824","2.210669660176497E-05","2.210669660176497E-05"),
        @("ECs","Ccl11","Ackr4","FAPs","2","0.6666666666666666","0.259924","0.779772","0.0006491957374851489","0.0006491957374851488","3","1","1.890486333333333","5.671459","0.9120863940138783","0.9120863940138783","0.4913827697053333","4.422444927348001","0.0005921225992120098","0.0005921225992120097"),
        @("ECs","Ccl11","Ackr4","MuSCs","2","0.6666666666666666","0.259924","0.779772","0.0006491957374851489","0.0006491957374851488","2","0.6666666666666666","0.1116383333333333","0.334915","0.05386116952465989","0.05386116952465989","0.02901748215333333","0.26115733938","3.49664416713742E-05","3.496644167137419E-05"),
        @("FAPs","Ccl11","Ackr4","ECs","3","1","383.1307676666667","1149.392303","0.9569214896224009","0.9569214896224006","1","0.3333333333333333","0.07058066666666667","0.211742","0.03405243646146196","0.03405243646146196","27.04162500242511","243.374625021826","0.03258550822397433","0.03258550822397433"),
        @("FAPs","Ccl11","Ackr4","FAPs","3","1","383.1307676666667","1149.392303","0.9569214896224009","0.9569214896224006","3","1","1.890486333333333","5.671459","0.9120863940138783","0.9120863940138783","724.3034801533421","6518.731321380078","0.8727950708240845","0.8727950708240843"),
        @("FAPs","Ccl11","Ackr4","MuSCs","3","1","383.1307676666667","1149.392303","0.9569214896224009","0.9569214896224006","2","0.6666666666666666","0.1116383333333333","0.334915","0.05386116952465989","0.05386116952465989","42.77208035102723","384.948723159245","0.0515409105743422","0.05154091057434219"),
        @("MuSCs","Ccl11","Ackr4","ECs","3","1","16.98778433333333","50.963353","0.042429314640114","0.04242931464011399","1","0.3333333333333333","0.07058066666666667","0.211742","0.03405243646146196","0.03405243646146196","1.199009143436222","10.791082290926","0.00144482154088586","0.00144482154088586"),
        @("MuSCs","Ccl11","Ackr4","FAPs","3","1","16.98778433333333","50.963353","0.042429314640114","0.04242931464011399","3","1","1.890486333333333","5.671459","0.9120863940138783","0.9120863940138783","32.11517411578078","289.036567042027","0.03869920059058183","0.03869920059058182"),
        @("MuSCs","Ccl11","Ackr4","MuSCs","3","1","16.98778433333333","50.963353","0.042429314640114","0.04242931464011399","2","0.6666666666666666","0.1116383333333333","0.334915","0.05386116952465989","0.05386116952465989","1.896487929999445","17.068391369995","0.002285292508646314","0.002285292508646313")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $col = 5
    for ($i = 4; $i -lt $rec.Length; $i++) {
        $ws.Cells.Item($row, $col).Value = [double]$rec[$i]
        $col++
    }
    $row++
}
